# Added Configurable zero_before_threshold parameter to enable setting dims
# before noise_threshold or First Rise Point to 0.
#
# This updates the computed First_Noticeable_Increase_Index (C),
# First_Noticeable_Increase_Cumulative_Value (E) and Pulse_Width (G)
# columns on each of the Step3_DataPts_* sheets to reflect the new
# zero_before_threshold behavior.

$wb = $excel.ActiveWorkbook

# Per-row updates that are identical across every Step3_DataPts_* sheet
# (columns C and E do not depend on the intensity threshold in column B).
$rowUpdates = @{
    2 = @{ C = 88; E = 0.01269704761990285 }
    3 = @{ C = 87; E = 0.01256054747815522 }
    4 = @{ C = 87; E = 0.01237489128571141 }
    5 = @{ C = 87; E = 0.01096804624834805 }
    6 = @{ C = 87; E = 0.01231409369821012 }
}

# Per-sheet, per-row Pulse_Width (column G) values, which do depend on the
# threshold used for each sheet.
$sheetGUpdates = @{
    "Step3_DataPts_0.5" = @{ 2 = 15; 3 = 23; 4 = 11; 5 = 23; 6 = 12 }
    "Step3_DataPts_0.7" = @{ 2 = 46; 3 = 52; 4 = 40; 5 = 53; 6 = 44 }
    "Step3_DataPts_0.8" = @{ 2 = 58; 3 = 61; 4 = 59; 5 = 66; 6 = 60 }
    "Step3_DataPts_0.9" = @{ 2 = 79; 3 = 79; 4 = 78; 5 = 79; 6 = 78 }
}

foreach ($sheetName in $sheetGUpdates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $rowUpdates.Keys) {
        $ws.Range("C$row").Value = $rowUpdates[$row].C
        $ws.Range("E$row").Value = $rowUpdates[$row].E
    }

    $gValues = $sheetGUpdates[$sheetName]
    foreach ($row in $gValues.Keys) {
        $ws.Range("G$row").Value = $gValues[$row]
    }
}
